$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix casing of the constructor names used in the "Example" column so the
# documentation is consistent with the lower-case type names used elsewhere
# in the sheet (cell, range, array, table, class).
$ws.Range("C2").Value = "let A1 = new cell(5);"
$ws.Range("C5").Value = "let tab1 = new table();"
